$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the header formatting (style index) from an existing header
# cell ("sum" in G1) onto the new H1 header cell by copying formats only,
# then set the text/values.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("H1").Value = "Save"

$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
